# Update market/profit data cells (H-N) across several sheets
# as refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 858.875
$ws.Range("I18").Value = 867.2857
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 867.2857
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = -583.2857
$ws.Range("N18").Value = -1368

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 719.4
$ws.Range("I32").Value = 629
$ws.Range("J32").Value = 779.6667
$ws.Range("K32").Value = 629
$ws.Range("L32").Value = 779.6667
$ws.Range("M32").Value = -303
$ws.Range("N32").Value = -1431.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2102.1052
$ws.Range("I40").Value = 1907.2727
$ws.Range("J40").Value = 2370
$ws.Range("K40").Value = 1907.2727
$ws.Range("L40").Value = 2370
$ws.Range("M40").Value = -1732.2727
$ws.Range("N40").Value = -2720

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1360.5
$ws.Range("I43").Value = 1561.1111
$ws.Range("J43").Value = 1240.1333
$ws.Range("K43").Value = 1561.1111
$ws.Range("L43").Value = 1240.1333
$ws.Range("M43").Value = -1492.1111
$ws.Range("N43").Value = -1378.1333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2822.6775
$ws.Range("I51").Value = 2143
$ws.Range("J51").Value = 3020.9167
$ws.Range("K51").Value = 2143
$ws.Range("L51").Value = 3020.9167
$ws.Range("M51").Value = -1659
$ws.Range("N51").Value = -3988.9167

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1941.0769
$ws.Range("I70").Value = 6175
$ws.Range("J70").Value = 1457.2
$ws.Range("K70").Value = 18525
$ws.Range("L70").Value = 4371.6
$ws.Range("M70").Value = -18255
$ws.Range("N70").Value = -4911.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1941.0769
$ws.Range("I73").Value = 6175
$ws.Range("J73").Value = 1457.2
$ws.Range("K73").Value = 18525
$ws.Range("L73").Value = 4371.6
$ws.Range("M73").Value = -17589
$ws.Range("N73").Value = -6243.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2890.9
$ws.Range("I74").Value = 2618.1667
$ws.Range("J74").Value = 3300
$ws.Range("K74").Value = 2618.1667
$ws.Range("L74").Value = 3300
$ws.Range("M74").Value = -1682.1667
$ws.Range("N74").Value = -5172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2890.9
$ws.Range("I77").Value = 2618.1667
$ws.Range("J77").Value = 3300
$ws.Range("K77").Value = 13090.8335
$ws.Range("L77").Value = 16500
$ws.Range("M77").Value = -8410.833500000001
$ws.Range("N77").Value = -25860

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 21024.666
$ws.Range("I125").Value = 25022.4
$ws.Range("J125").Value = 1036
$ws.Range("K125").Value = 225201.6
$ws.Range("L125").Value = 9324
$ws.Range("M125").Value = -222741.6
$ws.Range("N125").Value = -14244

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1542.85
$ws.Range("I135").Value = 607
$ws.Range("K135").Value = 5463
$ws.Range("M135").Value = -2928

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4459.189
$ws.Range("I137").Value = 1265.1666
$ws.Range("J137").Value = 5992.32
$ws.Range("K137").Value = 3795.4998
$ws.Range("L137").Value = 17976.96
$ws.Range("M137").Value = -1245.4998
$ws.Range("N137").Value = -23076.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 16
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4049.8333
$ws.Range("I86").Value = 3574.75
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3574.75
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2451.75
$ws.Range("N86").Value = -7246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4049.8333
$ws.Range("I89").Value = 3574.75
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 17873.75
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -12257.75
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 518.5
$ws.Range("I22").Value = 765.3333
$ws.Range("J22").Value = 271.66666
$ws.Range("K22").Value = 765.3333
$ws.Range("L22").Value = 271.66666
$ws.Range("M22").Value = -415.3333
$ws.Range("N22").Value = -971.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23283622
$ws.Range("I31").Value = 50001310
$ws.Range("J31").Value = 50845.824
$ws.Range("K31").Value = 50001310
$ws.Range("L31").Value = 50845.824
$ws.Range("M31").Value = -50001015
$ws.Range("N31").Value = -51435.824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23283622
$ws.Range("I34").Value = 50001310
$ws.Range("J34").Value = 50845.824
$ws.Range("K34").Value = 50001310
$ws.Range("L34").Value = 50845.824
$ws.Range("M34").Value = -50001108
$ws.Range("N34").Value = -51249.824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4012.5
$ws.Range("I62").Value = 4113.5137
$ws.Range("J62").Value = 2766.6667
$ws.Range("K62").Value = 4113.5137
$ws.Range("L62").Value = 2766.6667
$ws.Range("M62").Value = -3489.5137
$ws.Range("N62").Value = -4014.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4012.5
$ws.Range("I65").Value = 4113.5137
$ws.Range("J65").Value = 2766.6667
$ws.Range("K65").Value = 20567.5685
$ws.Range("L65").Value = 13833.3335
$ws.Range("M65").Value = -17447.5685
$ws.Range("N65").Value = -20073.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 615.7632
$ws.Range("I5").Value = 413.5238
$ws.Range("J5").Value = 865.58826
$ws.Range("K5").Value = 1240.5714
$ws.Range("L5").Value = 2596.76478
$ws.Range("M5").Value = -1128.5714
$ws.Range("N5").Value = -2820.76478

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1729.6666
$ws.Range("I122").Value = 690
$ws.Range("K122").Value = 6210
$ws.Range("M122").Value = -3760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 615.7632
$ws.Range("I135").Value = 413.5238
$ws.Range("J135").Value = 865.58826
$ws.Range("K135").Value = 3721.7142
$ws.Range("L135").Value = 7790.29434
$ws.Range("M135").Value = -1186.7142
$ws.Range("N135").Value = -12860.29434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2970.6667
$ws.Range("I22").Value = 1355.2
$ws.Range("J22").Value = 4990
$ws.Range("K22").Value = 1355.2
$ws.Range("L22").Value = 4990
$ws.Range("M22").Value = -1060.2
$ws.Range("N22").Value = -5580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2970.6667
$ws.Range("I27").Value = 1355.2
$ws.Range("J27").Value = 4990
$ws.Range("K27").Value = 1355.2
$ws.Range("L27").Value = 4990
$ws.Range("M27").Value = -1248.2
$ws.Range("N27").Value = -5204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2344.3333
$ws.Range("I46").Value = 5300.5
$ws.Range("J46").Value = 1499.7142
$ws.Range("K46").Value = 5300.5
$ws.Range("L46").Value = 1499.7142
$ws.Range("M46").Value = -5112.5
$ws.Range("N46").Value = -1875.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 500
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1727.5186
$ws.Range("I68").Value = 1589.6666
$ws.Range("J68").Value = 2210
$ws.Range("K68").Value = 1589.6666
$ws.Range("L68").Value = 2210
$ws.Range("M68").Value = -840.6666
$ws.Range("N68").Value = -3708

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1727.5186
$ws.Range("I71").Value = 1589.6666
$ws.Range("J71").Value = 2210
$ws.Range("K71").Value = 7948.333000000001
$ws.Range("L71").Value = 11050
$ws.Range("M71").Value = -4204.333000000001
$ws.Range("N71").Value = -18538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 48000
$ws.Range("J68").Value = 48000
$ws.Range("L68").Value = 48000
$ws.Range("N68").Value = -49622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 48000
$ws.Range("J71").Value = 48000
$ws.Range("L71").Value = 144000
$ws.Range("N71").Value = -152112
